$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (D=Fecha serial, J=Volumen, K=Precio minimo, L=Precio maximo, M=Precio promedio ponderado, P=Precio $/Kg)
$rows = @(
    @{ Row = 2;  D = 44984; J = 200; K = 17000; L = 18000; M = 17500; P = 972 },
    @{ Row = 3;  D = 44557; J = 400; K = 13000; L = 14000; M = 13500; P = 750 },
    @{ Row = 4;  D = 44957; J = 400; K = 21000; L = 22000; M = 21500; P = 1194 },
    @{ Row = 5;  D = 44964; J = 300; K = 20000; L = 21000; M = 20500; P = 1139 },
    @{ Row = 6;  D = 44977; J = 400; K = 16500; L = 17000; M = 16750; P = 931 },
    @{ Row = 7;  D = 45068; J = 400; K = 16000; L = 17000; M = 16500; P = 917 },
    @{ Row = 8;  D = 44547; J = 200; K = 13000; L = 14000; M = 13500; P = 750 },
    @{ Row = 9;  D = 44568; J = 500; K = 15000; L = 16000; M = 15500; P = 861 },
    @{ Row = 10; D = 44998; J = 320; K = 17000; L = 18000; M = 17500; P = 972 },
    @{ Row = 11; D = 44960; J = 400; K = 19500; L = 20000; M = 19750; P = 1097 },
    @{ Row = 12; D = 45005; J = 200; K = 17000; L = 18000; M = 17500; P = 972 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value = $r.D   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $r.J  # J - Volumen
    $ws.Cells.Item($row, 11).Value = $r.K  # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $r.L  # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $r.M  # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $r.P  # P - Precio $/Kg
}

$wb.Save()
